$d = $word.ActiveDocument

$replacements = @(
    @("2024-12-29 Sunday", "2024-12-30 Monday"),
    @("70×55=", "76×23="),
    @("78×50=", "46×55="),
    @("59×36=", "20×76="),
    @("22×18=", "76×77="),
    @("79×32=", "76×26="),
    @("73×87=", "36×74="),
    @("12×31=", "79×40="),
    @("25×14=", "13×77="),
    @("48×28=", "25×56="),
    @("46×77=", "12×47="),
    @("64×36=", "83×86="),
    @("63×49=", "18×60="),
    @("90×54=", "27×66="),
    @("88×55=", "27×38="),
    @("60×84=", "42×31="),
    @("14×90=", "33×72="),
    @("93×60=", "70×98="),
    @("27×78=", "39×90="),
    @("23×68=", "64×78="),
    @("19×66=", "68×17="),
    @("49×40=", "64×72="),
    @("33×73=", "85×83="),
    @("86×46=", "90×63="),
    @("97×41=", "51×63="),
    @("93×90=", "98×54=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
